$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "88-70="
$t.Cell(1,2).Range.Text = "5+81="
$t.Cell(1,3).Range.Text = "9+10="
$t.Cell(1,4).Range.Text = "88+8="
$t.Cell(1,5).Range.Text = "67-65="
$t.Cell(2,1).Range.Text = "17+75="
$t.Cell(2,2).Range.Text = "77-47="
$t.Cell(2,3).Range.Text = "3+88="
$t.Cell(2,4).Range.Text = "95-66="
$t.Cell(2,5).Range.Text = "23+4="
$t.Cell(3,1).Range.Text = "8+80="
$t.Cell(3,2).Range.Text = "38+18="
$t.Cell(3,3).Range.Text = "54-27="
$t.Cell(3,4).Range.Text = "92-30="
$t.Cell(3,5).Range.Text = "76-42="
$t.Cell(4,1).Range.Text = "2+78="
$t.Cell(4,2).Range.Text = "31+68="
$t.Cell(4,3).Range.Text = "53-11="
$t.Cell(4,4).Range.Text = "11+54="
$t.Cell(4,5).Range.Text = "1+80="
$t.Cell(5,1).Range.Text = "69+0="
$t.Cell(5,2).Range.Text = "76-37="
$t.Cell(5,3).Range.Text = "12+68="
$t.Cell(5,4).Range.Text = "41+31="
$t.Cell(5,5).Range.Text = "62-21="
$t.Cell(6,1).Range.Text = "39+3="
$t.Cell(6,2).Range.Text = "21+67="
$t.Cell(6,3).Range.Text = "31+37="
$t.Cell(6,4).Range.Text = "12+4="
$t.Cell(6,5).Range.Text = "47+32="
$t.Cell(7,1).Range.Text = "38+8="
$t.Cell(7,2).Range.Text = "25-11="
$t.Cell(7,3).Range.Text = "47+41="
$t.Cell(7,4).Range.Text = "91-41="
$t.Cell(7,5).Range.Text = "96-8="
$t.Cell(8,1).Range.Text = "33+54="
$t.Cell(8,2).Range.Text = "82-58="
$t.Cell(8,3).Range.Text = "49+11="
$t.Cell(8,4).Range.Text = "10+86="
$t.Cell(8,5).Range.Text = "71+19="
$t.Cell(9,1).Range.Text = "45+32="
$t.Cell(9,2).Range.Text = "10+53="
$t.Cell(9,3).Range.Text = "81-16="
$t.Cell(9,4).Range.Text = "98-41="
$t.Cell(9,5).Range.Text = "24+1="
$t.Cell(10,1).Range.Text = "6+85="
$t.Cell(10,2).Range.Text = "68-11="
$t.Cell(10,3).Range.Text = "48-19="
$t.Cell(10,4).Range.Text = "89-46="
$t.Cell(10,5).Range.Text = "66-32="
$t.Cell(11,1).Range.Text = "69-1="
$t.Cell(11,2).Range.Text = "12+87="
$t.Cell(11,3).Range.Text = "60+27="
$t.Cell(11,4).Range.Text = "28+71="
$t.Cell(11,5).Range.Text = "87-74="
$t.Cell(12,1).Range.Text = "61-6="
$t.Cell(12,2).Range.Text = "61+3="
$t.Cell(12,3).Range.Text = "81-45="
$t.Cell(12,4).Range.Text = "49-35="
$t.Cell(12,5).Range.Text = "74+5="
$t.Cell(13,1).Range.Text = "26-7="
$t.Cell(13,2).Range.Text = "39+46="
$t.Cell(13,3).Range.Text = "3+44="
$t.Cell(13,4).Range.Text = "56-0="
$t.Cell(13,5).Range.Text = "1+13="
$t.Cell(14,1).Range.Text = "67-60="
$t.Cell(14,2).Range.Text = "80-4="
$t.Cell(14,3).Range.Text = "63-41="
$t.Cell(14,4).Range.Text = "70-32="
$t.Cell(14,5).Range.Text = "82+11="
$t.Cell(15,1).Range.Text = "95-38="
$t.Cell(15,2).Range.Text = "77+10="
$t.Cell(15,3).Range.Text = "57-0="
$t.Cell(15,4).Range.Text = "6+6="
$t.Cell(15,5).Range.Text = "4+76="
$t.Cell(16,1).Range.Text = "1+87="
$t.Cell(16,2).Range.Text = "85+11="
$t.Cell(16,3).Range.Text = "23+5="
$t.Cell(16,4).Range.Text = "87-75="
$t.Cell(16,5).Range.Text = "55+14="
$t.Cell(17,1).Range.Text = "37+17="
$t.Cell(17,2).Range.Text = "96-5="
$t.Cell(17,3).Range.Text = "21+45="
$t.Cell(17,4).Range.Text = "23+28="
$t.Cell(17,5).Range.Text = "95-75="
$t.Cell(18,1).Range.Text = "10+44="
$t.Cell(18,2).Range.Text = "88-47="
$t.Cell(18,3).Range.Text = "17+32="
$t.Cell(18,4).Range.Text = "49-33="
$t.Cell(18,5).Range.Text = "42-21="
$t.Cell(19,1).Range.Text = "68-41="
$t.Cell(19,2).Range.Text = "57-54="
$t.Cell(19,3).Range.Text = "87-34="
$t.Cell(19,4).Range.Text = "7+60="
$t.Cell(19,5).Range.Text = "36+12="
$t.Cell(20,1).Range.Text = "87+1="
$t.Cell(20,2).Range.Text = "48+35="
$t.Cell(20,3).Range.Text = "47+49="
$t.Cell(20,4).Range.Text = "96-92="
$t.Cell(20,5).Range.Text = "58+37="
